$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the slightly corrected timestamp value in A13
$ws.Range("A13").Value = 45874.45852195602

# Add new row 14 with data, matching style of existing data rows (column A uses date/time style)
$ws.Range("A14").Value = 45874.50051116817
$ws.Range("B14").Value = 2025
$ws.Range("C14").Value = 19
$ws.Range("D14").Value = 18.17
$ws.Range("E14").Value = 80.14
$ws.Range("F14").Value = 652.15
$ws.Range("G14").Value = 12.61
$ws.Range("H14").Value = "ESE"
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = "12:00:44"

# Ensure A14 has the same style (number format) as the other date cells in column A
$ws.Range("A14").NumberFormat = $ws.Range("A13").NumberFormat
